# "fix error in linux"
# The Server.xlsx config sheet had a duplicate/bogus row (GameServer_2 on row 7,
# sharing GameServer_1's IP/name) and every server entry was still pointing at the
# loopback address 127.0.0.1. This cleans that up:
#   - give each real server row (2-6) its own LAN IP instead of 127.0.0.1
#   - wipe out the stray duplicate row 7
#   - autofit the IP column now that it holds longer addresses
#   - leave the selection sitting on the now-empty row 7 (where the cleanup happened)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-assign the IP column (F) for the five real server rows away from 127.0.0.1
# to distinct addresses.
$ws.Range("F2").Value = "192.168.1.113"
$ws.Range("F3").Value = "192.168.1.114"
$ws.Range("F4").Value = "192.168.1.115"
$ws.Range("F5").Value = "192.168.1.116"
$ws.Range("F6").Value = "192.168.1.117"

# Row 7 was a stray duplicate ("GameServer_2") left over from a copy/paste -
# clear it out entirely, keeping the existing cell formatting.
$ws.Range("A7:H7").ClearContents()

# The IP column now holds longer values ("192.168.1.11x") - autofit it so the
# column is wide enough to show them.
$ws.Columns("F:F").AutoFit()

# Leave the selection on the row that was just cleaned up.
$ws.Range("A7:XFD7").Select()
